# Update default max employee values and the (hidden) tables-under-plots rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "max employee" column (L) for rows 5-12, and let the
#     dependent "available hours" column (N) recalc automatically. ---
$ws.Range("L5").Value  = 200
$ws.Range("L6").Value  = 100
$ws.Range("L7").Value  = 100
$ws.Range("L8").Value  = 50
$ws.Range("L9").Value  = 50
$ws.Range("L10").Value = 50
$ws.Range("L11").Value = 50
$ws.Range("L12").Value = 50

# The old default (300) used a distinct highlighted/red-font "max" style
# (style index 6). The new defaults reuse the ordinary highlighted-input
# look already used elsewhere in the sheet (e.g. L3), so re-stamp the
# format for L5:L12 by copying it from L3.
$ws.Range("L3").Copy()
$ws.Range("L5:L12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add back the (blank, formatted) table row 13 under the plot area. ---
$ws.Range("L3").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the selected cell / viewport to reflect where the user was
#     last working. ---
$ws.Range("L14").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 9
